$wb = $excel.ActiveWorkbook

# ---------- Summary sheet ----------
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.9377224199288257
$ws1.Range("C2").Value = 0.4
$ws1.Range("D2").Value = 0.5
$ws1.Range("E2").Value = 0.4444444444444444
$ws1.Range("F2").Value = 0.4761904761904762
$ws1.Range("G2").Value = 0.4952380952380953
$ws1.Range("H2").Value = 0.7303370786516854
$ws1.Range("I2").Value = 14
$ws1.Range("J2").Value = 21
$ws1.Range("K2").Value = 513
$ws1.Range("L2").Value = 14

# ---------- Classification Report sheet ----------
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 ("0")
$ws2.Range("B2").Value = 0.9734345351043643
$ws2.Range("C2").Value = 0.9606741573033708
$ws2.Range("D2").Value = 0.9670122525918945

# row 3 ("1")
$ws2.Range("B3").Value = 0.4
$ws2.Range("C3").Value = 0.5
$ws2.Range("D3").Value = 0.4444444444444444

# row 4 ("accuracy")
$ws2.Range("B4").Value = 0.9377224199288257
$ws2.Range("C4").Value = 0.9377224199288257
$ws2.Range("D4").Value = 0.9377224199288257
$ws2.Range("E4").Value = 0.9377224199288257

# row 5 ("macro avg")
$ws2.Range("B5").Value = 0.6867172675521822
$ws2.Range("C5").Value = 0.7303370786516854
$ws2.Range("D5").Value = 0.7057283485181695

# row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.9448648429639334
$ws2.Range("C6").Value = 0.9377224199288257
$ws2.Range("D6").Value = 0.9409768457802777

# ---------- Confusion Matrix sheet ----------
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 513
$ws3.Range("C2").Value = 21

$ws3.Range("B3").Value = 14
$ws3.Range("C3").Value = 14
